$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace all "unfinished" placeholder text ("waffles") with "placeholder"
$ws.Range("F2").Value = "placeholder"
$ws.Range("F3").Value = "placeholder"
$ws.Range("F4").Value = "placeholder"

# Fill in the previously empty Diet cell for the Turkey Swiss Wrap row
$ws.Range("E3").Value = "NA"

# The underlying table definition referenced a trailing empty row (row 5) that
# has no data - shrink the table back down to the actual data range A1:F4
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:F4"))

# Update the active selection to match where the user last clicked
$ws.Range("D4").Select()
